$d = $word.ActiveDocument

# Remove the "_GoBack" bookmark from its current location
# (right before "14, intel i5" in the "On Macbook Pro mid-14, intel i5" paragraph).
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# Locate the spot between "Where one iteration" and " corresponds to one of
# the 3 RGB channel processing." in the final paragraph, and re-insert the
# "_GoBack" bookmark there (collapsed / zero-length range), splitting the
# run in two.
$r = $d.Content
$found = $r.Find.Execute("Where one iteration", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target = $d.Range($r.End, $r.End)
$d.Bookmarks.Add("_GoBack", $target)
